# Apply the "cosmic test cases, cravat_to_vcf" edit:
#  - insert a new column A (shifts existing A:D -> B:E)
#  - populate the new column A with group markers (1/2) on certain rows
#  - add "cosmic" test-case annotations to rows 39-43
#  - update the saved view (topLeftCell/zoom/selection)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift all existing columns one to the right by inserting a blank column A.
$ws.Columns("A:A").Insert()

# New column A markers (group numbers used by the test cases).
$ws.Range("A11").Value = 1
$ws.Range("A37").Value = 1
$ws.Range("A38").Value = 1
$ws.Range("A48").Value = 2
$ws.Range("A49").Value = 2
$ws.Range("A50").Value = 2
$ws.Range("A51").Value = 2
$ws.Range("A52").Value = 2
$ws.Range("A53").Value = 2
$ws.Range("A54").Value = 2
$ws.Range("A55").Value = 2
$ws.Range("A56").Value = 2

# New "cosmic" related test rows (39-43).
$ws.Range("C39").Value = "yes"
$ws.Range("D39").Value = "cosmic"

$ws.Range("E40").Value = "include in cosmic? database needs updating"
$ws.Range("E41").Value = "include in cosmic? database needs updating"

$ws.Range("C42").Value = "yes"

$ws.Range("E43").Value = "include in cosmic? database needs updating"

# Restore the view state: scroll position, zoom, and current selection.
$ws.Application.ActiveWindow.Zoom = 115
$ws.Application.ActiveWindow.ScrollRow = 25
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("C46").Select()
